$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "custom accuracy": row 5 (B5:AH5) values are rounded down to 2 decimal places
$ws.Range("B5").Value = 13.79
$ws.Range("C5").Value = 10.32
$ws.Range("D5").Value = 0.98
$ws.Range("E5").Value = 30.3
$ws.Range("F5").Value = 24.55
$ws.Range("G5").Value = 10.51
$ws.Range("H5").Value = 40.84
$ws.Range("I5").Value = 16.79
$ws.Range("J5").Value = 7.65
$ws.Range("K5").Value = 10.77
$ws.Range("L5").Value = 12.16
$ws.Range("M5").Value = 12.98
$ws.Range("N5").Value = 3.52
$ws.Range("O5").Value = 10.9
$ws.Range("P5").Value = 15.34
$ws.Range("Q5").Value = 9.34
$ws.Range("R5").Value = 0.31
$ws.Range("S5").Value = 0.58
$ws.Range("T5").Value = 158.61
$ws.Range("U5").Value = 30.42
$ws.Range("V5").Value = 10.06
$ws.Range("W5").Value = 20.29
$ws.Range("X5").Value = 10.85
$ws.Range("Y5").Value = 1.45
$ws.Range("Z5").Value = 20.49
$ws.Range("AA5").Value = 8.89
$ws.Range("AB5").Value = 7.94
$ws.Range("AC5").Value = 9.29
$ws.Range("AD5").Value = 12.84
$ws.Range("AE5").Value = 0.47
$ws.Range("AF5").Value = 36.9
$ws.Range("AG5").Value = 5.59
$ws.Range("AH5").Value = 12.58

# "데이터 1000개": the dataset was trimmed, dropping the last sample row (row 6)
$ws.Rows.Item(6).Delete()
